$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column G ("K") values for rows 2-15 per regenerated strike/K calculation.
$values = @{
    2  = 1
    3  = 1
    4  = 1
    5  = 0
    6  = 0
    7  = 2
    8  = 1
    9  = 1
    10 = 0
    12 = 2
    13 = 1
    14 = 2
    15 = 1
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
